$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('土地')
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 15
$ws.Range("B2").Value = '南投縣南投市牛運堀段02670002地號'
$ws.Range("C2").Value = 325
$ws.Range("D2").Value = '全部'
$ws.Range("E2").Value = '廖述嘉'
$ws.Range("F2").Value = '79年05月04日'
$ws.Range("G2").Value = '共有物分割'
$ws.Range("H2").Value = '(超過五年）'
$ws.Range("I2").Value = 'land'
$ws.Range("J2").Value = 'normal'
$ws.Range("K2").Value = '2012-02-10'
$ws.Range("L2").Value = '盧秀燕'
$ws.Range("M2").Value = 869
$ws.Range("N2").Value = 'tmp61a71'
$ws.Range("O2").Value = 15
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 325
$ws.Range("2:2").Style = $ws.Range("3:3").Style

$ws = $wb.Worksheets.Item('汽車')
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 32
$ws.Range("B2").Value = 'HYUNDAI'
$ws.Range("C2").Value = 2497
$ws.Range("D2").Value = '盧秀燕'
$ws.Range("E2").Value = '99年02月06日'
$ws.Range("F2").Value = '買賣'
$ws.Range("G2").Value = 100000
$ws.Range("2:2").Style = $ws.Range("3:3").Style

$ws = $wb.Worksheets.Item('存款')
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 48
$ws.Range("B2").Value = '臺灣銀行群賢分行'
$ws.Range("C2").Value = '活期儲蓄存款'
$ws.Range("D2").Value = '新臺幣'
$ws.Range("E2").Value = '盧秀燕'
$ws.Range("F2").Value = 5353477
$ws.Range("2:2").Style = $ws.Range("3:3").Style

$ws = $wb.Worksheets.Item('保險')
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 93
$ws.Range("B2").Value = '南山人壽'
$ws.Range("C2").Value = '子女教育保險'
$ws.Range("D2").Value = '廖述嘉'
$ws.Range("E2").Value = '保險期間:951811718(22年)年繳保費71400'
$ws.Range("2:2").Style = $ws.Range("3:3").Style
